{"js": "// Splits two paragraphs in the \"LOM3018\" syllabus document by inserting\n// manual line breaks (<w:br/>) at specific points, without altering any\n// of the surrounding text:\n//\n//   1. In \"Programa resumido\", a break is inserted right after\n//      \"...No\u00e7\u00f5es b\u00e1sicas de Projetos em Engenharia.\" and before\n//      \"Em todos o conte\u00fado do curso...\".\n//\n//   2. In \"Bibliografia\", a break is inserted after each of the first\n//      four numbered references, turning the single run-on paragraph\n//      into five visually separated lines.\n//\n// Word's JS API represents a manual line break as the vertical-tab\n// character (U+000B) in run text; inserting it (rather than calling\n// range.insertBreak, which in this host only lands at paragraph\n// boundaries) reliably places a real <w:br/> exactly at the search\n// match's position.\n\nconst body = context.document.body;\n\nasync function insertBreakAfter(searchText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `insertBreakAfter: expected exactly 1 match for ${JSON.stringify(searchText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(\"\\v\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 1) \"Programa resumido\" paragraph: split right before \"Em todos\".\nawait insertBreakAfter(\"No\u00e7\u00f5es b\u00e1sicas de Projetos em Engenharia.\");\n\n// 2) \"Bibliografia\" paragraph: split after each of the first four entries.\nawait insertBreakAfter(\"Editora, 2010.\");\nawait insertBreakAfter(\"Editora, 2006.\");\nawait insertBreakAfter(\"7a.ed., 2008. \");\nawait insertBreakAfter(\"UFSCar, 1985.\");\n", "ps1": "# Splits two paragraphs in the \"LOM3018\" syllabus document by inserting\n# manual line breaks (a literal <w:br/>) at specific points, without\n# altering any of the surrounding text:\n#\n#   1. In \"Programa resumido\", a break is inserted right after\n#      \"...No\u00e7\u00f5es b\u00e1sicas de Projetos em Engenharia.\" and before\n#      \"Em todos o conte\u00fado do curso...\".\n#\n#   2. In \"Bibliografia\", a break is inserted after each of the first\n#      four numbered references, turning the single run-on paragraph\n#      into five visually separated lines.\n#\n# Word represents a manual line break ([char]11, a.k.a. vertical tab /\n# Chr(11)) inside Range.Text; using Find to locate the anchor text, then\n# collapsing the found range to its end and calling InsertAfter([char]11)\n# places a genuine <w:br/> exactly at that position.\n\n$d = $word.ActiveDocument\n\nfunction Insert-BreakAfter([string]$needle) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Text = $needle\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Insert-BreakAfter: text not found: $needle\"\n    }\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter([char]11)\n}\n\n# 1) \"Programa resumido\" paragraph: split right before \"Em todos\".\nInsert-BreakAfter \"No\u00e7\u00f5es b\u00e1sicas de Projetos em Engenharia.\"\n\n# 2) \"Bibliografia\" paragraph: split after each of the first four entries.\nInsert-BreakAfter \"Editora, 2010.\"\nInsert-BreakAfter \"Editora, 2006.\"\nInsert-BreakAfter \"7a.ed., 2008. \"\nInsert-BreakAfter \"UFSCar, 1985.\"\n"}
